$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = 94154
$ws.Range("B36").Value = 2
$ws.Range("C36").Value = 0
